# Generate Report for Handoff
# The "70a42f8e-..." file is in "Ready for handoff" status; regenerating the
# handoff report refreshes its "Latest Handoff Datetime" (column D, row 5) on
# both the zh-cn and de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D5").Value = "2016-02-24 09:28:28"
$wsDe.Range("D5").Value = "2016-02-24 09:28:41"
